# Improvements_to_MKPedals.xlsx — "Added shopping cart and checkout functionality"
#
# The meaningful change in this revision is a new row appended to the
# improvements table on Sheet1 (row 25): a page-wide ("ALL") task about
# sorting files into folders, plus its note about needing to edit hrefs.
# The sheet's active window view also moved (new scroll position/selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in the new table row (Page / Improvement / Notes / Done? / Solution?) ---
# Column D (Done?) on row 25 already holds its FALSE checkbox value and is untouched.
$ws.Range("A25").Value = "ALL"
$ws.Range("B25").Value = "Sort out files into folders"
$ws.Range("C25").Value = "Will require editing all the hrefs within all the files"

# The extra text wraps onto a second line, so the row grows to match the
# other two-line rows (e.g. row 24) in this sheet.
$ws.Rows.Item(25).RowHeight = 29

# --- Restore the window's scroll position / active selection ---
[void]$ws.Range("E24").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
